$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Data History")
$ws2 = $wb.Worksheets.Item("Sheet1")
Write-Host $ws1.Name
Write-Host $ws2.Name
$v = $ws1.Range("E24").Value
Write-Host "E24 value:"
Write-Host $v
Write-Host ("Row128 height: " + $ws2.Rows.Item(128).RowHeight.ToString())
